$d = $word.ActiveDocument

# Locate the anchor paragraph: "Set the static folder to be called 'public' ..."
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Set the static folder*") {
        $anchor = $p
        break
    }
}

# Insert a new list paragraph right after the anchor (it inherits the
# ListParagraph style + numPr from the anchor paragraph).
$anchor.Range.InsertParagraphAfter()
$p1 = $anchor.Next()
$p1.Range.Text = "Getting Cannot GET/ error.  Need to fix this by writing the GET request for the website to get the homepage."

# Append the second run's text onto the end of the same paragraph (before
# its paragraph mark).
$p1End = $p1.Range.End
$ins = $d.Range($p1End - 1, $p1End - 1)
$ins.InsertAfter("  The res.render in this app.get is currently trying to access the index.ejs file, which does not currently exist.")

# Insert the second new list paragraph.
$p1.Range.InsertParagraphAfter()
$p2 = $p1.Next()
$p2.Range.Text = "Written and app.get res.render for the toDoListTable page as well."

# Insert the third new list paragraph.
$p2.Range.InsertParagraphAfter()
$p3 = $p2.Next()
$p3.Range.Text = "NOW WANT TO WRITE MY PARTIALS, WRITE THE INDEX.EJS AND TODOLISTTABLE.EJS PAGE, AND WRITE THE POST AND GET REQUESTS TO WORK WITH THE DATA FROM THE FORM AND SHOW IT ON THE TABLE ON THE TODOLISTTABLE PAGE."

Write-Output "done"
